$wb = $excel.ActiveWorkbook

# Add the new "test" worksheet after the last existing sheet ("Pay bill")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "test"

# New columns / headers on "General Data" sheet for the find-transactions scenario
$general = $wb.Worksheets.Item("General Data")
$general.Range("D1").Value = "Transaction date from"
$general.Range("E1").Value = "Transaction date to"
$general.Columns.Item(4).ColumnWidth = 20.5
$general.Columns.Item(5).ColumnWidth = 20.5

# Selection / active sheet bookkeeping to match the recorded session state
$general.Activate()
$general.Range("H7:H8").Select()
